$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list update (prices + 1h volume %) per source diff.
# D-column price cells are plain-text in the source (inlineStr) even when
# the text looks like a pure number (e.g. "0.627", "59.31"); Excel would
# otherwise auto-coerce such an assignment to a numeric cell, so we force
# text format, assign, then drop back to the Normal style (no explicit
# number format left behind) to match the original unstyled text cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.366.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.981.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.47%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("E5").Value = "  -3.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.627"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.94%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.31"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -9.36%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -1.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0826"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.55%  "
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.16%  "
$ws.Range("E14").Value = "  -7.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.270.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.986.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.259.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.83%  "
$ws.Range("E20").Value = "  -5.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0867"
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = "  -3.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.37%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  -2.70%  "
$ws.Range("E26").Value = "  -3.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.131"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.01%  "
$ws.Range("E31").Value = "  -1.75%  "
$ws.Range("E32").Value = "  -2.07%  "
$ws.Range("E33").Value = "  -6.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0627"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  -5.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.70%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("E38").Value = "  -7.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.51%  "
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("E42").Value = "  -6.99%  "
$ws.Range("E43").Value = "  -4.34%  "
$ws.Range("E44").Value = "  -2.07%  "
$ws.Range("E45").Value = "  -4.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "92.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.371.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.12%  "
$ws.Range("E49").Value = "  -5.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.60%  "
